# Applies updated cosinor analysis results (CircaDB / CircadiPy simulation rerun)
# to rows 2-8 of the active worksheet, per commit "Make figures again to publication".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = [double]"25.04000000000048"
$ws.Range("H2").Value = [double]"0.0005155467964560456"
$ws.Range("I2").Value = [double]"0.0005155467964560456"
$ws.Range("L2").Value = [double]"40.68046926439931"
$ws.Range("M2").Value = "[14.500862453316941, 66.86007607548169]"
$ws.Range("N2").Value = [double]"0.00306910020338913"
$ws.Range("O2").Value = [double]"0.00306910020338913"
$ws.Range("P2").Value = [double]"1.13839493553504"
$ws.Range("Q2").Value = "[0.4842895582110396, 1.7925003128590395]"
$ws.Range("R2").Value = [double]"0.001045003025312852"
$ws.Range("S2").Value = [double]"0.001045003025312852"
$ws.Range("T2").Value = [double]"61.4179174512302"
$ws.Range("U2").Value = "[47.650226803879065, 75.18560809858133]"
$ws.Range("V2").Value = [double]"1.330646703934235e-11"
$ws.Range("W2").Value = [double]"1.330646703934235e-11"
$ws.Range("X2").Value = [double]"20.50322322322361"
$ws.Range("Y2").Value = [double]"17.8964564564568"
$ws.Range("Z2").Value = [double]"23.10998998999043"
$ws.Range("F3").Value = [double]"25.04000000000048"
$ws.Range("H3").Value = [double]"0.0002401000911699258"
$ws.Range("I3").Value = [double]"0.0002401000911699258"
$ws.Range("L3").Value = [double]"49.50827873837769"
$ws.Range("M3").Value = "[20.907565848414137, 78.10899162834124]"
$ws.Range("N3").Value = [double]"0.001104629060933116"
$ws.Range("O3").Value = [double]"0.001104629060933116"
$ws.Range("P3").Value = [double]"1.301921279866041"
$ws.Range("Q3").Value = "[0.6603948521059619, 1.9434477076261194]"
$ws.Range("R3").Value = [double]"0.0001775047010990072"
$ws.Range("S3").Value = [double]"0.0001775047010990072"
$ws.Range("T3").Value = [double]"57.13425255464279"
$ws.Range("U3").Value = "[41.21635401389119, 73.05215109539438]"
$ws.Range("V3").Value = [double]"4.66374427965377e-09"
$ws.Range("W3").Value = [double]"4.66374427965377e-09"
$ws.Range("X3").Value = [double]"19.85153153153191"
$ws.Range("Y3").Value = [double]"17.29489489489522"
$ws.Range("Z3").Value = [double]"22.4081681681686"
$ws.Range("F4").Value = [double]"25.04000000000048"
$ws.Range("H4").Value = [double]"0.0007300076776983744"
$ws.Range("I4").Value = [double]"0.0007300076776983744"
$ws.Range("L4").Value = [double]"46.8772302197163"
$ws.Range("M4").Value = "[18.34529652680436, 75.40916391262823]"
$ws.Range("N4").Value = [double]"0.001847978366893432"
$ws.Range("O4").Value = [double]"0.001847978366893432"
$ws.Range("P4").Value = [double]"1.125815985971117"
$ws.Range("Q4").Value = "[0.45913165908319353, 1.7925003128590413]"
$ws.Range("R4").Value = [double]"0.001416979416367647"
$ws.Range("S4").Value = [double]"0.001416979416367647"
$ws.Range("T4").Value = [double]"70.49746211848964"
$ws.Range("U4").Value = "[54.205698321936296, 86.78922591504298]"
$ws.Range("V4").Value = [double]"3.210876009518415e-11"
$ws.Range("W4").Value = [double]"3.210876009518415e-11"
$ws.Range("X4").Value = [double]"20.55335335335374"
$ws.Range("Y4").Value = [double]"17.8964564564568"
$ws.Range("Z4").Value = [double]"23.21025025025069"
$ws.Range("F5").Value = [double]"25.04000000000048"
$ws.Range("H5").Value = [double]"6.147843897297278e-06"
$ws.Range("I5").Value = [double]"6.147843897297278e-06"
$ws.Range("L5").Value = [double]"67.45898674403264"
$ws.Range("M5").Value = "[37.31502059812459, 97.6029528899407]"
$ws.Range("N5").Value = [double]"4.648024516273885e-05"
$ws.Range("O5").Value = [double]"4.648024516273885e-05"
$ws.Range("P5").Value = [double]"0.4339737599553466"
$ws.Range("Q5").Value = "[-0.06918422260157797, 0.9371317425122712]"
$ws.Range("R5").Value = [double]"0.08919929066938792"
$ws.Range("S5").Value = [double]"0.08919929066938792"
$ws.Range("T5").Value = [double]"78.56114330806682"
$ws.Range("U5").Value = "[61.83133644865161, 95.29095016748202]"
$ws.Range("V5").Value = [double]"2.892575068358383e-12"
$ws.Range("W5").Value = [double]"2.892575068358383e-12"
$ws.Range("X5").Value = [double]"23.31051051051096"
$ws.Range("Y5").Value = [double]"21.30530530530571"
$ws.Range("Z5").Value = [double]"25.3157157157162"
$ws.Range("F6").Value = [double]"22"
$ws.Range("H6").Value = [double]"0.000223608428781108"
$ws.Range("I6").Value = [double]"0.000223608428781108"
$ws.Range("J6").Value = [double]"0.1026666312020235"
$ws.Range("K6").Value = [double]"0.1026666312020235"
$ws.Range("L6").Value = [double]"49.60268624729384"
$ws.Range("M6").Value = "[18.37996826341586, 80.82540423117182]"
$ws.Range("N6").Value = [double]"0.002521882286491639"
$ws.Range("O6").Value = [double]"0.002521882286491639"
$ws.Range("P6").Value = [double]"-0.5031579825569237"
$ws.Range("Q6").Value = "[-1.1446844103170015, 0.13836844520315417]"
$ws.Range("R6").Value = [double]"0.1211825134688984"
$ws.Range("S6").Value = [double]"0.1211825134688984"
$ws.Range("T6").Value = [double]"61.59551887993514"
$ws.Range("U6").Value = "[45.448553988617306, 77.74248377125298]"
$ws.Range("V6").Value = [double]"1.000676652651578e-09"
$ws.Range("W6").Value = [double]"1.000676652651578e-09"
$ws.Range("X6").Value = [double]"1.761761761761761"
$ws.Range("Y6").Value = [double]"-0.4844844844844851"
$ws.Range("Z6").Value = [double]"4.008008008008007"
$ws.Range("F7").Value = [double]"22"
$ws.Range("H7").Value = [double]"0.0001097973048674872"
$ws.Range("I7").Value = [double]"0.0001097973048674872"
$ws.Range("J7").Value = [double]"0.9613698801638011"
$ws.Range("K7").Value = [double]"0.9613698801638011"
$ws.Range("L7").Value = [double]"43.29311225083914"
$ws.Range("M7").Value = "[21.795615858231358, 64.79060864344692]"
$ws.Range("N7").Value = [double]"0.0001957937965528789"
$ws.Range("O7").Value = [double]"0.0001957937965528789"
$ws.Range("P7").Value = [double]"0.1195000208572692"
$ws.Range("Q7").Value = "[-0.5220264069028095, 0.7610264486173479]"
$ws.Range("R7").Value = [double]"0.7092920626824832"
$ws.Range("S7").Value = [double]"0.7092920626824832"
$ws.Range("T7").Value = [double]"60.19844365315244"
$ws.Range("U7").Value = "[46.80456186714295, 73.59232543916193]"
$ws.Range("V7").Value = [double]"1.069011545951071e-11"
$ws.Range("W7").Value = [double]"1.069011545951071e-11"
$ws.Range("X7").Value = [double]"21.58158158158158"
$ws.Range("Y7").Value = [double]"19.33533533533533"
$ws.Range("Z7").Value = [double]"23.82782782782783"
$ws.Range("F8").Value = [double]"22"
$ws.Range("H8").Value = [double]"0.0203634592813825"
$ws.Range("I8").Value = [double]"0.0203634592813825"
$ws.Range("J8").Value = [double]"0.7089832252745661"
$ws.Range("K8").Value = [double]"0.7089832252745661"
$ws.Range("L8").Value = [double]"33.13681898397434"
$ws.Range("M8").Value = "[2.2179658559771696, 64.0556721119715]"
$ws.Range("N8").Value = [double]"0.03625365235520395"
$ws.Range("O8").Value = [double]"0.03625365235520395"
$ws.Range("P8").Value = [double]"0.4968685077749626"
$ws.Range("Q8").Value = "[-0.69184222601577, 1.6855792415656952]"
$ws.Range("R8").Value = [double]"0.4043105236025499"
$ws.Range("S8").Value = [double]"0.4043105236025499"
$ws.Range("T8").Value = [double]"56.99372984715013"
$ws.Range("U8").Value = "[40.67336068599528, 73.31409900830498]"
$ws.Range("V8").Value = [double]"9.085047381063305e-09"
$ws.Range("W8").Value = [double]"9.085047381063305e-09"
$ws.Range("X8").Value = [double]"20.26026026026026"
$ws.Range("Y8").Value = [double]"16.0980980980981"
$ws.Range("Z8").Value = [double]"24.42242242242242"

Write-Output "Applied all cell updates"
